$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "IIII" example entry to "IIII*III" (C4). Dependent formulas in
# C6 (_xlfn.CONCAT("1-",LEN(C4))) and D6 (LEN(C4)+1) recalc automatically.
$ws.Range("C4").Value = "IIII*III"

# Minor numeric corrections in row 4
$ws.Range("E4").Value = 0
$ws.Range("I4").Value = 0

# Move the active selection from D8 to C4
[void]$ws.Range("C4").Select()
